# Fine tuning of the pair-integrate actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Remove the "0% refactored, 0% erroneous code" row entirely (old row 52);
# everything below shifts up by one row.
$ws.Rows.Item(52).Delete()

# Tweak a few wording/values in the remaining "pair-integration" action table.
$ws.Range("A46").Value = "No coding standard, 0% refactored, 86% erroneous code (everything wrong)"
$ws.Range("A51").Value = "10% erroneous code"
$ws.Range("A54").Value = "53% refactored, 3% erroneous code"

# Fill in the Duration (B) / Errors introduced into code (C) columns for
# every action row in the pair-integration table.
$ws.Range("B44").Value = 24
$ws.Range("C44").Value = 4

$ws.Range("B45").Value = 17
$ws.Range("C45").Value = 3

$ws.Range("B46").Value = 34
$ws.Range("C46").Value = "> 12"
$ws.Range("C46").HorizontalAlignment = -4152

$ws.Range("B47").Value = 21
$ws.Range("C47").Value = 4

$ws.Range("B48").Value = 21
$ws.Range("C48").Value = 27

$ws.Range("B49").Value = 29
$ws.Range("C49").Value = 14

$ws.Range("B50").Value = 26
$ws.Range("C50").Value = "?"
$ws.Range("C50").HorizontalAlignment = -4152

$ws.Range("B51").Value = 18
$ws.Range("C51").Value = 12

$ws.Range("B52").Value = 38
$ws.Range("C52").Value = 28

$ws.Range("B53").Value = 18
$ws.Range("C53").Value = 8

$ws.Range("B54").Value = 19
$ws.Range("C54").Value = 16

# Restore scroll position / selection similar to the authored view state.
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B52").Select()
